# Add a new test-data row to the "ContactUs" sheet (row 7) and make
# that sheet the active / selected sheet with the new row selected,
# matching the authored commit "Added a row in dynamic data".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ContactUs")

$ws.Cells.Item(7, 1).Value = "John@123#"
$ws.Cells.Item(7, 2).Value = "validemail@example.com"
$ws.Cells.Item(7, 3).Value = "Invalid Name TC"
$ws.Cells.Item(7, 4).Value = "Testing invalid characters."
$ws.Cells.Item(7, 5).Value = "INVALID_NAME"

$ws.Activate()
$ws.Range("A7:E7").Select()
